$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (the 2009 year row); subsequent rows shift up.
$ws.Rows("2").Delete()
